$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: ticket price becomes unavailable ("不可售") instead of a numeric price
    $ws.Range("G2").Value = "不可售"

    # Row 3: interested-count increases
    $ws.Range("F3").Value = 2138

    # Row 4: interested-count increases
    $ws.Range("F4").Value = 1622

    # Row 7: interested-count increases, and ticket price becomes available again (65)
    $ws.Range("F7").Value = 516
    $ws.Range("G7").Value = 65

    # Row 8: interested-count increases
    $ws.Range("F8").Value = 27

    # Row 9: interested-count increases
    $ws.Range("F9").Value = 5680
}
